$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H62").Value = 235227.67
$ws.Range("I62").Value = 200005
$ws.Range("J62").Value = 252839
$ws.Range("K62").Value = 200005
$ws.Range("L62").Value = 252839
$ws.Range("M62").Value = -199381
$ws.Range("N62").Value = -254087

$ws.Range("H64").Value = 3323.8667
$ws.Range("I64").Value = 2780
$ws.Range("J64").Value = 3595.8
$ws.Range("K64").Value = 2780
$ws.Range("L64").Value = 3595.8
$ws.Range("M64").Value = -2532
$ws.Range("N64").Value = -4091.8

$ws.Range("H65").Value = 235227.67
$ws.Range("I65").Value = 200005
$ws.Range("J65").Value = 252839
$ws.Range("K65").Value = 1000025
$ws.Range("L65").Value = 1264195
$ws.Range("M65").Value = -996905
$ws.Range("N65").Value = -1270435

$ws.Range("H67").Value = 3323.8667
$ws.Range("I67").Value = 2780
$ws.Range("J67").Value = 3595.8
$ws.Range("K67").Value = 2780
$ws.Range("L67").Value = 3595.8
$ws.Range("M67").Value = -1922
$ws.Range("N67").Value = -5311.8

$ws.Range("H125").Value = 2339.3333
$ws.Range("I125").Value = 3000
$ws.Range("J125").Value = 2207.2
$ws.Range("K125").Value = 27000
$ws.Range("L125").Value = 19864.8
$ws.Range("M125").Value = -24540
$ws.Range("N125").Value = -24784.8

$ws.Range("H134").Value = 46125
$ws.Range("J134").Value = 46125
$ws.Range("L134").Value = 46125
$ws.Range("N134").Value = -56265

$ws.Range("H141").Value = 3588
$ws.Range("I141").Value = 3681.2778
$ws.Range("J141").Value = 3028.3333
$ws.Range("K141").Value = 11043.8334
$ws.Range("L141").Value = 9084.999899999999
$ws.Range("M141").Value = -5863.8334
$ws.Range("N141").Value = -19444.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5301.4
$ws.Range("I102").Value = 4891
$ws.Range("J102").Value = 6122.2
$ws.Range("K102").Value = 4891
$ws.Range("L102").Value = 6122.2
$ws.Range("M102").Value = -3269
$ws.Range("N102").Value = -9366.200000000001

$ws.Range("H128").Value = 44143.5
$ws.Range("J128").Value = 44143.5
$ws.Range("L128").Value = 44143.5
$ws.Range("N128").Value = -54103.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9180.4
$ws.Range("J50").Value = 9180.4
$ws.Range("L50").Value = 9180.4
$ws.Range("N50").Value = -10430.4

$ws.Range("H59").Value = 11663.75
$ws.Range("I59").Value = 9000
$ws.Range("J59").Value = 12551.667
$ws.Range("K59").Value = 9000
$ws.Range("L59").Value = 12551.667
$ws.Range("M59").Value = -7855
$ws.Range("N59").Value = -14841.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 512.0833
$ws.Range("I5").Value = 512.0833
$ws.Range("K5").Value = 1536.2499
$ws.Range("M5").Value = -1424.2499

$ws.Range("H113").Value = 920.433
$ws.Range("I113").Value = 568.7222
$ws.Range("J113").Value = 1000.56964
$ws.Range("K113").Value = 1706.1666
$ws.Range("L113").Value = 3001.70892
$ws.Range("M113").Value = 463.8334
$ws.Range("N113").Value = -7341.70892

$ws.Range("H117").Value = 144771.42
$ws.Range("I117").Value = 2000
$ws.Range("J117").Value = 201880
$ws.Range("K117").Value = 6000
$ws.Range("L117").Value = 605640
$ws.Range("M117").Value = -2558
$ws.Range("N117").Value = -612524

$ws.Range("H122").Value = 542.5833
$ws.Range("I122").Value = 279.22223
$ws.Range("J122").Value = 1332.6666
$ws.Range("K122").Value = 2513.00007
$ws.Range("L122").Value = 11993.9994
$ws.Range("M122").Value = -63.00007000000005
$ws.Range("N122").Value = -16893.9994

$ws.Range("H129").Value = 23811246
$ws.Range("J129").Value = 33334276
$ws.Range("L129").Value = 100002828
$ws.Range("N129").Value = -100012828

$ws.Range("H135").Value = 512.0833
$ws.Range("I135").Value = 512.0833
$ws.Range("K135").Value = 4608.7497
$ws.Range("M135").Value = -2073.7497

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 111556890
$ws.Range("J33").Value = 111556890
$ws.Range("L33").Value = 111556890
$ws.Range("N33").Value = -111557394

$ws.Range("H126").Value = 23814724
$ws.Range("I126").Value = 5600
$ws.Range("J126").Value = 55560224
$ws.Range("K126").Value = 16800
$ws.Range("L126").Value = 166680672
$ws.Range("M126").Value = -14330
$ws.Range("N126").Value = -166685612

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3144.476
$ws.Range("I7").Value = 2211.7778
$ws.Range("J7").Value = 3844
$ws.Range("K7").Value = 2211.7778
$ws.Range("L7").Value = 3844
$ws.Range("M7").Value = -2099.7778
$ws.Range("N7").Value = -4068

$ws.Range("H34").Value = 2250
$ws.Range("I34").Value = 500
$ws.Range("K34").Value = 500
$ws.Range("M34").Value = -328

$ws.Range("H40").Value = 3364.5
$ws.Range("J40").Value = 3729
$ws.Range("L40").Value = 3729
$ws.Range("N40").Value = -4001

$ws.Range("H122").Value = 2608.9
$ws.Range("I122").Value = 1750
$ws.Range("J122").Value = 3181.5
$ws.Range("K122").Value = 5250
$ws.Range("L122").Value = 9544.5
$ws.Range("M122").Value = -2800
$ws.Range("N122").Value = -14444.5

$ws.Range("H126").Value = 3144.476
$ws.Range("I126").Value = 2211.7778
$ws.Range("J126").Value = 3844
$ws.Range("K126").Value = 6635.3334
$ws.Range("L126").Value = 11532
$ws.Range("M126").Value = -6635.3334
$ws.Range("N126").Value = -16472

$ws.Range("H132").Value = 4130.5
$ws.Range("I132").Value = 4236.0884
$ws.Range("K132").Value = 12708.2652
$ws.Range("M132").Value = -10178.2652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6299
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 6299
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 6299
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -7339

$ws.Range("H126").Value = 2155.8
$ws.Range("I126").Value = 1910.1111
$ws.Range("J126").Value = 2524.3333
$ws.Range("K126").Value = 5730.3333
$ws.Range("L126").Value = 7572.999899999999
$ws.Range("M126").Value = -3260.3333
$ws.Range("N126").Value = -12512.9999
